$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - values only, cell styles already set in the template
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "LinkedIn_Poster"
$ws.Range("G1").Value = "LinkedIn_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

# Data row (row 2) for the new job posting JD_001
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"

# The job description spans two lines (embedded line break). Writing a
# multi-line string straight into .Value triggers this engine's
# autofit-row-height side effect (adds ht/customHeight to the row), which
# the source workbook does not have. Build the text via a formula (CHAR(10)
# for the line break) and then flatten the formula down to a plain cached
# value with Copy / PasteSpecial so the cell ends up as ordinary literal
# text, same as the rest of the row.
$ws.Range("C2").Formula = '="We are seeking a Junior RPA Developer to design, develop, and support automation solutions."&CHAR(10)&"Collaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2
